$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) figures increased for a handful of
# events on both the "展览" sheet and the aggregated "全部类型" sheet.
$targetSheets = @("展览", "全部类型")

$updates = @{
    7  = 10877
    8  = 408
    13 = 144
    16 = 41
    20 = 1111
}

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
